# Update benchmark: 2026-01-09 06:43:26 UTC
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (ŞANS OYUNLARI)
$ws.Range("E2").Value = ""
$ws.Range("K2").Value = "32,62 TL - 32,62 TL"

# Row 3 (HESAPTAN EFT - Şube)
$ws.Range("F3").Value = ""

# Row 4 (HESAPTAN EFT - ATM)
$ws.Range("F4").Value = ""

# Row 5 (HESAPTAN EFT - Mobil)
$ws.Range("F5").Value = ""

# Row 8 (HESAPTAN HAVALE - Şube)
$ws.Range("F8").Value = ""

# Row 9 (HESAPTAN HAVALE - ATM)
$ws.Range("F9").Value = ""

# Row 10 (HESAPTAN HAVALE - Mobil)
$ws.Range("F10").Value = ""

# Row 13 (GELEN SWIFT)
$ws.Range("C13").Value = "Hesaba: Asgari 0 TL | Azami 9.999.999.999.999 TL"
$ws.Range("D13").Value = "Hesaba: Asgari 1 TL | Azami 300 TL"
$ws.Range("F13").Value = ""

# Row 14 (GİDEN SWIFT - Mobil)
$ws.Range("F14").Value = ""

# Row 15 (ÇEK TAHSİLİ BAŞKA BANKA)
$ws.Range("K15").Value = "%0,3 Asgari Tutar: 237,26 TL Azami Tutar: 237,26 TL / 298,96 TL"

# Row 17 (AYNI ŞUBE ÇEK TAHSİLATI)
$ws.Range("K17").Value = "%0,6 Asgari Tutar: 237,26 TL Azami Tutar: 237,26 TL / 3.034,67 TL"

# Row 20 (ÇEK İADE)
$ws.Range("K20").Value = "147,11 TL"

# Row 21 (BLOKE ÇEK DÜZENLEME)
$ws.Range("K21").Value = "%0,9 Asgari Tutar: 446,06 TL Azami Tutar: 446,06 TL / 2.427,26 TL"

# Row 22 (YP ÇEK TAKASA GÖNDERME)
$ws.Range("K22").Value = "%0,3 Asgari Tutar: 73,56 TL Azami Tutar: 73,56 TL / 9.115,86 TL"

# Row 23 (ÇEK KARNESİ SAYFA ÜCRETİ)
$ws.Range("K23").Value = "64,8 TL"

# Row 24 (SENET TAHSİLE ALMA)
$ws.Range("K24").Value = "446,06 TL"

# Row 25 (MUAMELESİZ SENET İADESİ)
$ws.Range("K25").Value = "374,4 TL"
